# Adds a new "Meningococcal" row (with MenACWY / MenB special-situations
# detail) to the "adults" sheet, plus the "See Conditions and Alternate
# Dosing" label referenced by it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adults")

# ---- Row 13: main Meningococcal data row -----------------------------
$ws.Range("A13").Value = "Meningococcal "
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "Conditional"
$ws.Range("D13").Value = "See Conditions and Alternate Dosing"
$ws.Range("E13").Value = "See Conditions and Alternate Dosing"
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Range("H13").Value = "X"
$ws.Range("I13").Value = 6935
$ws.Range("J13").Value = 364635
$ws.Range("K13").Value = 6935
$ws.Range("L13").Value = 364635
$ws.Range("M13").Value = 6935
$ws.Range("N13").Value = 364635

$ws.Range("U13").Value = "for MenB ->  Anatomical or functional asplenia, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis. For MenACWY -> Anatomical or functional asplenia, HIV infection, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis, Travel to countries with high rates of meningococcal disease, First-year college students living in residential housing, Military"
$ws.Range("U13").Value = "for MenB -> Anatomical or functional asplenia, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis. For MenACWY -> Anatomical or functional asplenia, HIV infection, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis, Travel to countries with high rates of meningococcal disease, First-year college students living in residential housing, Military"

# Rich-text cell: "Pregnancy" (normal) + ": " (bold) + rest (normal)
$ws.Range("V13").Value = "Pregnancy: Delay MenB until after pregnancy unless at increased risk and vaccination benefits outweigh potential risks "
$ws.Range("V13").Font.Name = "Helvetica"
$ws.Range("V13").Font.Size = 11
$ws.Range("V13").Font.Color = 1973537
$ws.Range("V13").Characters(10, 2).Font.Bold = $true

$ws.Range("W13").Value = "Special situations for MenACWY:"

$ws.Range("Y13").Value = "Disease states, Anatomical or functional asplenia, HIV infection, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis"
$ws.Range("Z13").Value = "2 doses at least 8 weeks apart, every 5 years"
$ws.Range("AA13").Value = "Travel to countries with high rates of meningococcal disease"
$ws.Range("AB13").Value = "1 dose every 5 years"
$ws.Range("AC13").Value = "First-year college students living in residential housing, Military"
$ws.Range("AD13").Value = "1 dose"

$ws.Range("AE13").Value = "Special situations for MenB:"
$ws.Range("AG13").Value = "Anatomical or functional asplenia, Persistent complement component deficiency, Complement inhibitor (ex: eculizumab, ravulizumab) use, Microbiologist exposed to Neisseria meningitidis"
$ws.Range("AH13").Value = "Option 1: Bexsero (MenB-4C) -> 2 doses with 1 month spacing, Booster: 1 year after primary series, repeated every 2-3 years if risk remains "

# ---- Row 14: continuation of MenB detail ------------------------------
$ws.Range("AH14").Value = "------- Option 2: Trumenba (MenB-FHbp) -> 2-4 doses -> Dose 2 with 1-2 months spacing, Dose 3 is not needed if dose 2 was > 6 months after dose 1; if needed it should be 4 months after dose 2. Dose 4: Only needed if dose 3 was < 4 months after dose 2; should be > 4 months after dose 3, Booster: 1 year after primary series, repeated every 2-3 years if risk remains"

# ---- Formatting: blank cells in AF:AT for rows 13-15 carry the plain
# Calibri 11 black style (matching the sheet's other "condition" columns)
$blankFont = $ws.Range("V3").Font  # existing cell that already uses this look

foreach ($row in 13..15) {
    $cols = @("AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AR", "AS", "AT")
    foreach ($col in $cols) {
        $addr = "$col$row"
        if ($ws.Range($addr).Value -eq $null -or $ws.Range($addr).Value -eq "") {
            $ws.Range($addr).Font.Name = "Calibri"
            $ws.Range($addr).Font.Size = 11
            $ws.Range($addr).Font.Color = 0
        }
    }
}

# ---- Sheet bookkeeping -------------------------------------------------
$ws.Range("T15").Select()

Write-Output "done"
